$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Update first table (rows 2-9), columns D (Estimated SP Sprint 2) and E (Estimated SP Sprint 3), F (Estimated SP Sprint 4) ---
# Row 3 - Lider de desarrollo
$ws.Range("D3").Value = 48
$ws.Range("E3").Value = 8
$ws.Range("F3").Value = 45

# Row 4 - Dev y adm de BDD
$ws.Range("D4").Value = 69
$ws.Range("E4").Value = 16

# Row 5 - Dev full stack A
$ws.Range("D5").Value = 88
$ws.Range("E5").Value = 40

# Row 7 - Tester A
$ws.Range("D7").Value = 70
$ws.Range("E7").Value = 16

# --- Update second table (rows 20-27), columns E (Estimated SP Sprint 2), G (Estimated SP Sprint 3), I (Estimated SP Sprint 4) ---
# Row 21 - Lider de desarrollo
$ws.Range("E21").Value = 48
$ws.Range("G21").Value = 8
$ws.Range("I21").Value = 45

# Row 22 - Dev y adm de BDD
$ws.Range("E22").Value = 69
$ws.Range("G22").Value = 16

# Row 23 - Dev full stack A
$ws.Range("E23").Value = 88
$ws.Range("G23").Value = 40

# Row 25 - Tester A
$ws.Range("E25").Value = 70
$ws.Range("G25").Value = 16

# --- Selection on the active sheet view ---
$ws.Range("G14").Select()

$wb.Save()
